# New Sale method: getTotalPayed(Sale)
#
# Inserts a new API-reference row right after the "closeSale" row (row 30),
# in the blank separator row that used to be row 31, pushing every row
# below it down by one. The new row documents:
#   RETURN TYPE       number(,2)
#   METHOD SIGNATURE  getSaleTotalPayed(SaleRef)
#   VERB              GET
#   URL               /sale/{saleId}/payed   (as a formula, like its siblings)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 31 (and everything after it) down by one row, Excel-style.
$ws.Rows(31).EntireRow.Insert()

$ws.Range("A31").Value = "number(,2)"
$ws.Range("B31").Value = "getSaleTotalPayed(SaleRef)"
$ws.Range("C31").Value = "GET"
$ws.Range("D31").Formula = '="/sale/{saleId}/payed"'

# Column D in this table carries the "green" banding fill (same as the
# other GET rows in the Sale block, e.g. D28); copy that formatting onto
# the new D31 cell instead of inventing a brand-new style.
$ws.Range("D28").Copy()
$ws.Range("D31").PasteSpecial(-4122)
